# iDVC-FreeDIC application form — update the "邮箱" (Email) header cell so
# applicants are reminded to use their institutional email address.
#
# D2 currently reads "邮箱". We change it to "邮箱（请使用机构邮箱）" where the
# inserted reminder "请使用机构邮箱" is emphasised in red, matching the rest
# of the sheet's plain black text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "邮箱（请使用机构邮箱）"

# "请使用机构邮箱" starts at character 4 (1-based) and is 7 characters long.
# 255 == RGB(255, 0, 0) as a VBA/COM color long (red).
$cell.Characters(4, 7).Font.Color = 255

# Leave the trailing "）" in the sheet's normal (automatic/black) colour.
$cell.Characters(11, 1).Font.ColorIndex = 1

# Reflect the cursor position left behind at the end of the editing session.
$ws.Range("E13").Select() | Out-Null
